# Update cryptos list figures (price/volume) per latest scrape, and fix
# row 12/13 and 18/19 coin ordering (TRON/WrappedEther, Avalanche/ShibaInu swap).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '30.978.18'
$ws.Range("E2").Value = '  -0.30%  '
$ws.Range("D3").Value = "'" + '1.955.07'
$ws.Range("E3").Value = '  -0.67%  '
$ws.Range("D4").Value = "'" + '1.001'
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").Value = "'" + '243.80'
$ws.Range("E5").Value = '  -2.06%  '
$ws.Range("D6").Value = "'" + '1.001'
$ws.Range("E6").Value = '  +0.26%  '
$ws.Range("D7").Value = "'" + '0.4834'
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("E8").Value = '  -0.30%  '
$ws.Range("D9").Value = "'" + '0.07045'
$ws.Range("E9").Value = '  +3.39%  '
$ws.Range("D10").Value = "'" + '19.74'
$ws.Range("E10").Value = '  +2.22%  '
$ws.Range("D11").Value = "'" + '107.18'
$ws.Range("E11").Value = '  -2.34%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = "'" + '1.974.95'
$ws.Range("E12").Value = '  +0.59%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = "'" + '0.07786'
$ws.Range("E13").Value = '  +0.44%  '
$ws.Range("D14").Value = "'" + '5.440'
$ws.Range("E14").Value = '  -1.01%  '
$ws.Range("D15").Value = "'" + '0.6998'
$ws.Range("E15").Value = '  +0.26%  '
$ws.Range("D16").Value = "'" + '280.10'
$ws.Range("E16").Value = '  -3.91%  '
$ws.Range("D17").Value = "'" + '30.984.21'
$ws.Range("E17").Value = '  -0.28%  '
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = "'" + '0.000007806'
$ws.Range("E18").Value = '  +0.74%  '
$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D19").Value = "'" + '13.27'
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("D20").Value = "'" + '2.198.99'
$ws.Range("E20").Value = '  -0.96%  '
$ws.Range("E21").Value = '  +0.12%  '
$ws.Range("D22").Value = "'" + '5.538'
$ws.Range("E22").Value = '  -2.29%  '
$ws.Range("D23").Value = "'" + '1.002'
$ws.Range("E23").Value = '  +0.26%  '
$ws.Range("D24").Value = "'" + '6.491'
$ws.Range("E24").Value = '  -2.24%  '
$ws.Range("D25").Value = "'" + '9.805'
$ws.Range("E25").Value = '  -1.21%  '
$ws.Range("D26").Value = "'" + '168.91'
$ws.Range("E26").Value = '  -0.70%  '
$ws.Range("D27").Value = "'" + '19.75'
$ws.Range("E27").Value = '  -1.99%  '
$ws.Range("D28").Value = "'" + '2.170'
$ws.Range("E28").Value = '  -0.78%  '
$ws.Range("D29").Value = "'" + '0.1046'
$ws.Range("E30").Value = '  -4.17%  '
$ws.Range("E31").Value = '  -2.06%  '
$ws.Range("D32").Value = "'" + '4.601'
$ws.Range("E32").Value = '  -3.70%  '
$ws.Range("D33").Value = "'" + '4.439'
$ws.Range("E33").Value = '  -1.20%  '
$ws.Range("D34").Value = "'" + '0.04895'
$ws.Range("E34").Value = '  -4.06%  '
$ws.Range("D35").Value = "'" + '0.7458'
$ws.Range("E35").Value = '  -4.40%  '
$ws.Range("D36").Value = "'" + '1.163'
$ws.Range("E36").Value = '  -1.37%  '
$ws.Range("D37").Value = "'" + '2.736'
$ws.Range("E37").Value = '  +0.33%  '
$ws.Range("D38").Value = "'" + '0.01996'
$ws.Range("E38").Value = '  -2.49%  '
$ws.Range("D39").Value = "'" + '2.688'
$ws.Range("E39").Value = '  -1.52%  '
$ws.Range("D40").Value = "'" + '6.506'
$ws.Range("E40").Value = '  -0.59%  '
$ws.Range("D41").Value = "'" + '77.68'
$ws.Range("E41").Value = '  +8.52%  '
$ws.Range("D42").Value = "'" + '2.119'
$ws.Range("E42").Value = '  -0.90%  '
$ws.Range("E43").Value = '  +1.48%  '
$ws.Range("D44").Value = "'" + '109.03'
$ws.Range("E44").Value = '  -1.23%  '
$ws.Range("D45").Value = "'" + '0.4437'
$ws.Range("E45").Value = '  -0.69%  '
$ws.Range("D46").Value = "'" + '7.993'
$ws.Range("E46").Value = '  +6.11%  '
$ws.Range("E47").Value = '  +0.15%  '
$ws.Range("D48").Value = "'" + '990.46'
$ws.Range("E48").Value = '  +6.14%  '
$ws.Range("D49").Value = "'" + '9.355'
$ws.Range("E49").Value = '  -0.52%  '
$ws.Range("E50").Value = '  -2.70%  '
$ws.Range("D51").Value = "'" + '35.85'
$ws.Range("E51").Value = '  -0.49%  '
